# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage (no numeric
# auto-conversion), then restore the cell to its original (default/General)
# style so no stray number-format style is left behind.
function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "68.257.70"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "3.920.66"
$ws.Range("E3").Value = "  -0.66%  "

Set-TextValue $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.06%  "

Set-TextValue $ws "D5" "487.28"
$ws.Range("E5").Value = "  +3.35%  "

Set-TextValue $ws "D6" "147.59"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("E7").Value = "  -0.75%  "

$ws.Range("E8").Value = "  -0.04%  "

Set-TextValue $ws "D9" "0.734"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").Value = "  +1.87%  "

Set-TextValue $ws "D11" "0.0000350"
$ws.Range("E11").Value = "  +4.13%  "

Set-TextValue $ws "D12" "43.08"
$ws.Range("E12").Value = "  -0.58%  "

Set-TextValue $ws "D13" "10.78"
$ws.Range("E13").Value = "  +4.17%  "

$ws.Range("D14").Value = "4.551.42"
$ws.Range("E14").Value = "  -0.46%  "

$ws.Range("D15").Value = "3.944.98"
$ws.Range("E15").Value = "  +0.15%  "

Set-TextValue $ws "D16" "14.43"
$ws.Range("E16").Value = "  -5.65%  "

$ws.Range("E17").Value = "  -0.73%  "

Set-TextValue $ws "D18" "19.95"
$ws.Range("E18").Value = "  +0.44%  "

Set-TextValue $ws "D19" "1.14"
$ws.Range("E19").Value = "  -2.11%  "

$ws.Range("D20").Value = "68.423.63"
$ws.Range("E20").Value = "  +1.38%  "

Set-TextValue $ws "D21" "442.55"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D22" "15.23"
$ws.Range("E22").Value = "  +4.69%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws "D23" "3.49"
$ws.Range("E23").Value = "  +2.30%  "

Set-TextValue $ws "D24" "88.15"
$ws.Range("E24").Value = "  +0.46%  "

Set-TextValue $ws "D25" "11.48"
$ws.Range("E25").Value = "  +18.09%  "

Set-TextValue $ws "D26" "11.55"
$ws.Range("E26").Value = "  +13.45%  "

$ws.Range("E27").Value = "  -0.13%  "

Set-TextValue $ws "D28" "38.63"
$ws.Range("E28").Value = "  -1.07%  "

Set-TextValue $ws "D29" "5.74"
$ws.Range("E29").Value = "  -0.61%  "

Set-TextValue $ws "D30" "716.16"
$ws.Range("E30").Value = "  -0.87%  "

Set-TextValue $ws "D31" "13.84"
$ws.Range("E31").Value = "  +2.15%  "

$ws.Range("E32").Value = "  -0.55%  "

$ws.Range("E33").Value = "  +2.66%  "

Set-TextValue $ws "D34" "6.33"
$ws.Range("E34").Value = "  +18.24%  "

Set-TextValue $ws "D35" "42.26"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("D36").Value = "0.0₃0879"
$ws.Range("E36").Value = "  +11.81%  "

Set-TextValue $ws "D37" "61.41"
$ws.Range("E37").Value = "  +6.09%  "

Set-TextValue $ws "D38" "0.416"
$ws.Range("E38").Value = "  +23.77%  "

$ws.Range("E39").Value = "  -2.36%  "

$ws.Range("E40").Value = "  +19.49%  "

$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D41" "1.00"
$ws.Range("E41").Value = "  +0.07%  "

$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws "D42" "3.29"
$ws.Range("E42").Value = "  +7.76%  "

$ws.Range("E43").Value = "  +0.72%  "

Set-TextValue $ws "D44" "2.92"
$ws.Range("E44").Value = "  +4.60%  "

$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("E46").Value = "  +0.14%  "

Set-TextValue $ws "D47" "3.33"
$ws.Range("E47").Value = "  +6.80%  "

Set-TextValue $ws "D48" "3.44"
$ws.Range("E48").Value = "  -1.68%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D49" "2.14"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0350"
$ws.Range("E50").Value = "  +28.68%  "

Set-TextValue $ws "D51" "146.01"
$ws.Range("E51").Value = "  -0.30%  "

Write-Output "Updated $(($ws.UsedRange.Rows.Count)) rows of crypto data."
